# Update Gnas-Gcgr LR-pairs sheet with newly computed TPM-based NATMI values.
# The "Target cluster" set expands from {FAPs, MuSCs} to {ECs, FAPs, MuSCs}
# for every Sending cluster, so the whole 3x3 cross-join of
# Sending cluster x Target cluster (ECs/FAPs/MuSCs) is rewritten with the
# refreshed metric columns E:T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#           E..T (Ligand/Receptor/Edge metrics)
$data = @(
    @("ECs",   "Gnas", "Gcgr", "ECs",   3, 1, 68.158272,  204.474816,         0.164824640128582,  0.1648246401285819, 3, 1, 0.193739,            0.5812170000000001, 0.178467861985615,  0.178467861985615,  13.204915459008,    118.844239131072,   0.02941590112629642, 0.02941590112629642),
    @("ECs",   "Gnas", "Gcgr", "FAPs",  3, 1, 68.158272,  204.474816,         0.164824640128582,  0.1648246401285819, 3, 1, 0.8229729999999998, 2.468919,            0.7581035918523759, 0.7581035918523759, 56.09241758265598,  504.8317582439039,  0.1249541517072532,  0.1249541517072532),
    @("ECs",   "Gnas", "Gcgr", "MuSCs", 3, 1, 68.158272,  204.474816,         0.164824640128582,  0.1648246401285819, 1, 0.3333333333333333, 0.068856,           0.206568,            0.0634285461620092, 0.0634285461620092, 4.693105976831999,  42.237953791488,    0.01045458729503232, 0.01045458729503231),
    @("FAPs",  "Gnas", "Gcgr", "ECs",   3, 1, 143.539174, 430.617522,         0.3471155005059974, 0.3471155005059974, 3, 1, 0.193739,            0.5812170000000001, 0.178467861985615,  0.178467861985615,  27.809136031586,    250.282224284274,   0.06194896123737201, 0.061948961237372),
    @("FAPs",  "Gnas", "Gcgr", "FAPs",  3, 1, 143.539174, 430.617522,         0.3471155005059974, 0.3471155005059974, 3, 1, 0.8229729999999998, 2.468919,            0.7581035918523759, 0.7581035918523759, 118.128864644302,   1063.159781798718,  0.2631495077212319,  0.2631495077212318),
    @("FAPs",  "Gnas", "Gcgr", "MuSCs", 3, 1, 143.539174, 430.617522,         0.3471155005059974, 0.3471155005059974, 1, 0.3333333333333333, 0.068856,           0.206568,            0.0634285461620092, 0.0634285461620092, 9.883533364944,     88.951800284496,    0.02201703154739358, 0.02201703154739358),
    @("MuSCs", "Gnas", "Gcgr", "ECs",   3, 1, 201.822474, 605.4674219999999, 0.4880598593654206, 0.4880598593654206, 3, 1, 0.193739,            0.5812170000000001, 0.178467861985615,  0.178467861985615,  39.100884290286,    351.907958612574,   0.08710299962194654, 0.08710299962194652),
    @("MuSCs", "Gnas", "Gcgr", "FAPs",  3, 1, 201.822474, 605.4674219999999, 0.4880598593654206, 0.4880598593654206, 3, 1, 0.8229729999999998, 2.468919,            0.7581035918523759, 0.7581035918523759, 166.0944468952019,  1494.850022056818,  0.3699999324238908,  0.3699999324238907),
    @("MuSCs", "Gnas", "Gcgr", "MuSCs", 3, 1, 201.822474, 605.4674219999999, 0.4880598593654206, 0.4880598593654206, 1, 0.3333333333333333, 0.068856,           0.206568,            0.0634285461620092, 0.0634285461620092, 13.896688269744,    125.070194427696,   0.0309569273195833,  0.0309569273195833)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value2 = $values[$col - 1]
    }
}
